$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New use case: a single-IP "Source" (row 5, previously an empty gap row)
# paired with a multi-line "Destination" list.
$ws.Range("C5").Value = "9.208.46.246`n9.208.48.81`n9.208.48.48"
$ws.Range("A5").Value = "9.214.16.167"
$ws.Range("C5").WrapText = $true
$ws.Rows(5).RowHeight = 51

# Existing row 6 (156.78.90.3) gets its Destination value replaced with a
# comma-separated multi-line variant of the same IP list.
$ws.Range("C6").Value = "9.208.46.246,9.208.48.81,`n9.208.48.48"
$ws.Range("C6").WrapText = $true
$ws.Rows(6).RowHeight = 51

# Existing row 7 ((abcd fg)) gets an empty, wrap-styled Destination cell to
# match the new formatting used above.
$ws.Range("C7").WrapText = $true

# New row 8: duplicate the Source value from row 6 as a new use case entry.
$ws.Range("A8").Value = "156.78.90.3"

# Update the visible selection to span A3:A4.
$ws.Range("A3:A4").Select()

# Best-effort: restore the window position metadata (not guaranteed to be
# persisted by every host, but harmless to attempt).
$win = $excel.ActiveWindow
$win.Left = 27040
$win.Top = -20940
